$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the new row 4 first (sCs -> Ccl2/Ackr2 -> FAPs, recalculated values)
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ccl2"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.205968666666667
$ws.Range("H4").Value = 27.617906
$ws.Range("I4").Value = 0.1584804277009834
$ws.Range("J4").Value = 0.1584804277009834
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.51723133333333
$ws.Range("N4").Value = 34.551694
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 106.0272707814182
$ws.Range("R4").Value = 954.245437032764
$ws.Range("S4").Value = 0.1584804277009834
$ws.Range("T4").Value = 0.1584804277009834

# Overwrite row 3 in place (FAPs -> Ccl2/Ackr2 -> FAPs, recalculated values)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl2"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 25.22419466666667
$ws.Range("H3").Value = 75.672584
$ws.Range("I3").Value = 0.4342336264580882
$ws.Range("J3").Value = 0.4342336264580882
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.51723133333333
$ws.Range("N3").Value = 34.551694
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 290.5128851730329
$ws.Range("R3").Value = 2614.615966557296
$ws.Range("S3").Value = 0.4342336264580882
$ws.Range("T3").Value = 0.4342336264580882

# Overwrite row 2 in place (ECs -> Ccl2/Ackr2 -> FAPs, new sending cluster)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl2"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 23.658831
$ws.Range("H2").Value = 70.976493
$ws.Range("I2").Value = 0.4072859458409285
$ws.Range("J2").Value = 0.4072859458409285
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.51723133333333
$ws.Range("N2").Value = 34.551694
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 272.484229703238
$ws.Range("R2").Value = 2452.358067329142
$ws.Range("S2").Value = 0.4072859458409285
$ws.Range("T2").Value = 0.4072859458409285
